$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.930.92"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.39"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.31"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6948"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07691"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3046"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.39"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07817"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.13"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.831.03"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.099"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6815"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.584"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008261"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.933.12"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.25"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.074.67"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.460"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1507"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.34"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.756"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.538"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.214"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.175"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7824"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.853"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.144"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.298.70"
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01860"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9546"
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.157"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.07"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.684"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5164"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.975.82"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.88"
$ws.Range("E49").Value = "  -6.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.753"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.961"
$ws.Range("E51").Value = "  -1.14%  "
